$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated / new data rows 2-17 (Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster, then the 16 numeric NATMI metric columns E:T)
$data = @(
    @("ECs", "Mst1", "Mst1r", "ECs", 3, 1, 0.5054393333333334, 1.516318, 0.1786582562510795, 0.1786582562510796, 2, 0.6666666666666666, 0.19157, 0.5747100000000001, 0.04039874491219805, 0.04039874491219806, 0.09682701308666668, 0.8714431177800002, 0.007217569320745475, 0.007217569320745478),
    @("ECs", "Mst1", "Mst1r", "FAPs", 3, 1, 0.5054393333333334, 1.516318, 0.1786582562510795, 0.1786582562510796, 3, 1, 3.688005666666667, 11.064017, 0.777735554431318, 0.777735554431318, 1.864063125489556, 16.776568129406, 0.1389488779791658, 0.1389488779791659),
    @("ECs", "Mst1", "Mst1r", "M2", 3, 1, 0.5054393333333334, 1.516318, 0.1786582562510795, 0.1786582562510796, 1, 0.3333333333333333, 0.06405566666666666, 0.192167, 0.01350821390534768, 0.01350821390534768, 0.03237625345622222, 0.291386281106, 0.002413353941396001, 0.002413353941396002),
    @("ECs", "Mst1", "Mst1r", "sCs", 3, 1, 0.5054393333333334, 1.516318, 0.1786582562510795, 0.1786582562510796, 3, 1, 0.7983476666666668, 2.395043, 0.1683574867511363, 0.1683574867511363, 0.4035163124082223, 3.631646811674, 0.03007845500977225, 0.03007845500977225),
    @("FAPs", "Mst1", "Mst1r", "ECs", 3, 1, 1.665224, 4.995672000000001, 0.5886087537853821, 0.5886087537853822, 2, 0.6666666666666666, 0.19157, 0.5747100000000001, 0.04039874491219805, 0.04039874491219806, 0.3190069616800001, 2.871062655120001, 0.02377905489726244, 0.02377905489726245),
    @("FAPs", "Mst1", "Mst1r", "FAPs", 3, 1, 1.665224, 4.995672000000001, 0.5886087537853821, 0.5886087537853822, 3, 1, 3.688005666666667, 11.064017, 0.777735554431318, 0.777735554431318, 6.141355548269335, 55.27219993442401, 0.4577819554684013, 0.4577819554684014),
    @("FAPs", "Mst1", "Mst1r", "M2", 3, 1, 1.665224, 4.995672000000001, 0.5886087537853821, 0.5886087537853822, 1, 0.3333333333333333, 0.06405566666666666, 0.192167, 0.01350821390534768, 0.01350821390534768, 0.1066670334693333, 0.9600033012240001, 0.007951052952693065, 0.007951052952693069),
    @("FAPs", "Mst1", "Mst1r", "sCs", 3, 1, 1.665224, 4.995672000000001, 0.5886087537853821, 0.5886087537853822, 3, 1, 0.7983476666666668, 2.395043, 0.1683574867511363, 0.1683574867511363, 1.329427694877334, 11.964849253896, 0.09909669046702535, 0.09909669046702536),
    @("M2", "Mst1", "Mst1r", "ECs", 3, 1, 0.4572946666666667, 1.371884, 0.1616405023344417, 0.1616405023344417, 2, 0.6666666666666666, 0.19157, 0.5747100000000001, 0.04039874491219805, 0.04039874491219806, 0.08760393929333335, 0.7884354536400001, 0.006530073421288664, 0.006530073421288667),
    @("M2", "Mst1", "Mst1r", "FAPs", 3, 1, 0.4572946666666667, 1.371884, 0.1616405023344417, 0.1616405023344417, 3, 1, 3.688005666666667, 11.064017, 0.777735554431318, 0.777735554431318, 1.686505322003111, 15.178547898028, 0.1257135657016338, 0.1257135657016338),
    @("M2", "Mst1", "Mst1r", "M2", 3, 1, 0.4572946666666667, 1.371884, 0.1616405023344417, 0.1616405023344417, 1, 0.3333333333333333, 0.06405566666666666, 0.192167, 0.01350821390534768, 0.01350821390534768, 0.02929231473644444, 0.263630832628, 0.002183474481301489, 0.00218347448130149),
    @("M2", "Mst1", "Mst1r", "sCs", 3, 1, 0.4572946666666667, 1.371884, 0.1616405023344417, 0.1616405023344417, 3, 1, 0.7983476666666668, 2.395043, 0.1683574867511363, 0.1683574867511363, 0.3650801301124445, 3.285721171012001, 0.02721338873021779, 0.0272133887302178),
    @("sCs", "Mst1", "Mst1r", "ECs", 3, 1, 0.2011266666666667, 0.60338, 0.07109248762909652, 0.07109248762909653, 2, 0.6666666666666666, 0.19157, 0.5747100000000001, 0.04039874491219805, 0.04039874491219806, 0.03852983553333334, 0.3467685198000001, 0.002872047272901466, 0.002872047272901467),
    @("sCs", "Mst1", "Mst1r", "FAPs", 3, 1, 0.2011266666666667, 0.60338, 0.07109248762909652, 0.07109248762909653, 3, 1, 3.688005666666667, 11.064017, 0.777735554431318, 0.777735554431318, 0.7417562863844445, 6.67580657746, 0.055291155282117, 0.05529115528211701),
    @("sCs", "Mst1", "Mst1r", "M2", 3, 1, 0.2011266666666667, 0.60338, 0.07109248762909652, 0.07109248762909653, 1, 0.3333333333333333, 0.06405566666666666, 0.192167, 0.01350821390534768, 0.01350821390534768, 0.01288330271777778, 0.11594972446, 0.0009603325299571194, 0.0009603325299571196),
    @("sCs", "Mst1", "Mst1r", "sCs", 3, 1, 0.2011266666666667, 0.60338, 0.07109248762909652, 0.07109248762909653, 3, 1, 0.7983476666666668, 2.395043, 0.1683574867511363, 0.1683574867511363, 0.1605690050377778, 1.44512104534, 0.01196895254412094, 0.01196895254412094)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

Write-Output "Updated $($data.Length) rows"
